# "diem danh ngay 24" - add attendance column for 24/03/2012 (serial 40992) in column H
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header in H1 (Saturday 24/3/2012), same style/format as the other date cells (G1 etc.)
$ws.Range("H1").Value = 40992

# Mark attendance ("1") for each student row, mirroring column G's pattern
$ws.Range("H2").Value = "1"
$ws.Range("H3").Value = "1"
$ws.Range("H4").Value = "1"
$ws.Range("H5").Value = "1"
$ws.Range("H6").Value = "1"
$ws.Range("H7").Value = "1"

# Extend column H to the same display width as the preceding date columns (C:G)
$ws.Range("H1").EntireColumn.ColumnWidth = 9.83

# Move the active selection to H7, matching where the last edit was made
[void]$ws.Range("H7").Select()
